$wb = $excel.ActiveWorkbook

# Worksheet references
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 2 and 5
$wsOverview.Range("G2").Value = "2016-08-16 06:14:30"
$wsOverview.Range("G5").Value = "2016-08-16 06:14:30"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-16 06:14:25"
$wsZhCn.Range("H5").Value = "2016-08-16 06:14:25"
$wsZhCn.Range("K2").Value = "2016-08-16 06:14:50"
$wsZhCn.Range("K5").Value = "2016-08-16 06:14:50"

# de-de sheet: Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe.Range("H2").Value = "2016-08-16 06:14:30"
$wsDeDe.Range("H5").Value = "2016-08-16 06:14:30"
$wsDeDe.Range("K2").Value = "2016-08-16 06:14:56"
$wsDeDe.Range("K5").Value = "2016-08-16 06:14:56"
